$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update expected test results to reflect the fix where message bundle
# properties are correctly resolved instead of showing raw keys.

# Step1 result ( "= msg(\"say.hello\", 1, 2, 3)" ) now resolves
$ws.Range("D25:F25").Value = "Hello, from Project!"

# Step2 result ( "= msg(\"say.hello.1\", \"Parameter\")" ) now resolves
$ws.Range("D26:F26").Value = "Hello, Parameter!"

# Step3 result for the French locale column now resolves to the French text
$ws.Range("E27").Value = "Bonjour, from MessageBundle!"

# Step4 result for the French locale column now resolves to the French text
$ws.Range("E28").Value = "Bonjour, Parameter!"

# Step5 result ( "= msg(\"jar.say.hello\", null)" ) now resolves
$ws.Range("D29:F29").Value = "Hello, from Jar!"

# Step6 result ( "= msg(\"jar.say.hello.1\", \"Parameter\")" ) now resolves
$ws.Range("D30:F30").Value = "Hello, Parameter!"

# Column D widened so the resolved messages are fully visible
$ws.Columns.Item(4).ColumnWidth = 33

# Restore the view state (scroll position / active selection) as left by the editor
$ws.Range("E17").Select()
$excel.ActiveWindow.ScrollRow = 10
